$d = $word.ActiveDocument

# --- Paragraph 1: "Blazor Webassembly - Perfect for SPA and much more" ---
# Wrap "Blazor" and "Webassembly" with spell-check proofErr markers (split into runs).
$xmlP1 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Blazor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Webassembly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Perfect for SPA and much more</w:t></w:r></w:p></pkg:xmlData>'
$d.Paragraphs.Item(1).Range.InsertXML($xmlP1)

# --- Paragraph 2: the three "- ..." sentences joined by line breaks ---
$xmlP2 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>WebAssembly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a client-side technology that enables near-native performance without plug-ins in all modern browsers, including mobile.</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Blazor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a client-side library that uses .NET on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>WebAssembly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to support Single Page Applications written in C# using Razor templates.</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">- .NET 6 is coming and, with it, some notable </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Blazor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> improvements.</w:t></w:r></w:p></pkg:xmlData>
'@
$d.Paragraphs.Item(2).Range.InsertXML($xmlP2)

# --- Hyperlink paragraph for rId7 (optasy.com link) gains bold formatting ---
# Insert the bold runs first (InsertXML cannot retain rStyle), then re-apply the
# Hyperlink character style onto the link's own Range so w:rStyle survives alongside w:b/w:bCs.
$xmlP8 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:hyperlink r:id="rId7" w:history="1"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>https://www.optasy.com/blog/webassembly-vs-javascript-wasm-faster-js-when-does-javascript-perform-better</w:t></w:r></w:hyperlink><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></pkg:xmlData>
'@
$d.Paragraphs.Item(8).Range.InsertXML($xmlP8)
$hlP8 = $d.Hyperlinks.Item(4)
$hlP8.Range.Style = "Hyperlink"

# --- Paragraph 12: standalone "Blazor" ---
$xmlP12 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Blazor</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></pkg:xmlData>'
$d.Paragraphs.Item(12).Range.InsertXML($xmlP12)

# --- Paragraph 19: "Todo list" -> "Todo" (spellchecked) + " list" ---
$xmlP19 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Todo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> list</w:t></w:r></w:p></pkg:xmlData>'
$d.Paragraphs.Item(19).Range.InsertXML($xmlP19)

# --- Paragraph 34: "- Razor component file names require ... Todo.razor." ---
# Note: the source text uses non-breaking spaces (U+00A0) between several words;
# those are preserved verbatim below.
$xmlP34 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:t xml:space="preserve">- </w:t></w:r><w:r><w:t>Razor component file names require a capitalized first letter. Open the Pages folder and confirm that the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Todo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t> component file name starts with a capital letter T. The file name should be </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Todo.razor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></pkg:xmlData>
'@
$d.Paragraphs.Item(34).Range.InsertXML($xmlP34)

# --- Paragraph 37: "1.What is wasm" ---
$xmlP37 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:t xml:space="preserve">1.What is </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wasm</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></pkg:xmlData>'
$d.Paragraphs.Item(37).Range.InsertXML($xmlP37)

# --- Paragraph 38: "2.What is blazor" ---
$xmlP38 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:t xml:space="preserve">2.What is </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>blazor</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></pkg:xmlData>'
$d.Paragraphs.Item(38).Range.InsertXML($xmlP38)

# --- Paragraph 39: "3.Blazor webasm vs Blazor server" ---
$xmlP39 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:t xml:space="preserve">3.Blazor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>webasm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> vs </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Blazor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> server</w:t></w:r></w:p></pkg:xmlData>'
$d.Paragraphs.Item(39).Range.InsertXML($xmlP39)

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
